$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

$ws.Range("C31").Value = '1. In terms of the Code style, for each variable or function we defined, there are comments distinctly shows what it is dealing with. ' + $nl + '2. The image can be selected as a file and the system will read the picture and transform it into base64 and store the image and the News report information in MongoDB.' + $nl + '3. As for the annotation part, I used the socket to initiate the canvas and set the style display of the image to none to sync and create the same canvas in chatting rooms with the same number. And then I captured the data of the drawing action in the emit room and utilising socket to sync all canvases.' + $nl + '4. When a new people joins in, he will only receive the discussion information from the living socket. He will not see the previous annotations.' + $nl + '5. When I implemented the main page, I catched the data by axios post the service URL and captured the data returned from the controller. And then I used document create new tabs to encapsulate the data and dynamically add them to the news dashboard. At the same time I also created the onclick submit form so that users can click the picture or news title to create a new chatting room.' + $nl + '6. The service logic is that people can submit their new report to store in the MongoDB and browse all the news in the main page. When someone is interested in the news, he can on click the title or picture to start a chatting room. People in the same number Chatting room can have the discussion, searching information and taking annotations. ' + $nl + '7. People can preview all news reported on the main page to choose a news that they are intersted to start a chat.' + $nl + '8.The annotations created by the user, including the chatting history, and the drawn annotations on the picture, will be stored in indexedDB. These annotations will not be transferred to the database. Only when the user logs in again in the same room will they see their own annotations. In addition, indexedDB also contains stories and knowledge graphs published by users.' + $nl + '9.When offline, users can also upload stories, annotate stories, etc., and the data will be uploaded to indexedDB.'

$ws.Range("C42").Value = '1. We use axios to post the server and deal with the data returned in the .then function. When there is some error we use .catch function to capture the failure information and give out the alert.' + $nl + '2. In the new report page, users can fill in a form to submit their news and then it will be serialized and passed to the axios function as the data.' + $nl + '3. When there is new post to the newReport Server, it will call the controller to insert news into the database.' + $nl + '4. When the post request is sent to the main Page, the controller will be called to return the stored news to the client to display the data in the front-end pages.' + $nl + '5. We have tested the documents passing in the system, and we checked the data in MongoDB Compass, it can correctly stores the documents.' + $nl + '6. The routes configuration is well designed, for each get and post request it correctly render the page and when the post request is asking the service, it will also correctly pass some parameters and call the right functions.' + $nl + '7. When the server returns a failure with status code like 404 or 500, the client will catch the error information and gives out the alerts.' + $nl + '8. All the data passing through clients and servers are using JSON. We request the server using axios.post.' + $nl + '9. I also added the automatical data passing when user clicks a title ot image to start a chatting room. People can find the selected image url has been filled in the blank and it will make the system more convenient.'

$ws.Range("C51").Value = '1. Users can upload a new story by typing author, title, report content and select the picture file in the new report webpage. ' + $nl + '2. The information of the story will be stored into the MongoDB database. ' + $nl + '3. The story will be given an unique id automatically by MongoDB, which is used to get a specific story. ' + $nl + '4. The story will be given the detailed upload time automatically by MongoDB. ' + $nl + '5. The stories in the MongoDB database can be show on the main page, users can browse them in the page. ' + $nl + '6. The stories may be updated when users are browsing them, so we set a button to allow users refresh and update the news list.'

$ws.Range("C56").Value = '1. When user submits a new report form, ther serialized data will be passed to the axios, and the client will send post request to the "/newReport" url to call the controller to insert new information.' + $nl + '2. When people are browsing the main page, there will be a request to get news from the database. We used axios to post the server and get the data from controller.' + $nl + '3. The first utilisation of Socket.io is in chatting room, when people join in the chatting session everyone in the room will receive a note. And when someone sends a message, we use socket.io to add on everybody''s history panel.' + $nl + '4. Another Utilisation of socket.io is to sync the canvas. I initiate the canvas and hide the origin image in the chatting rooms. And when someone is drawing the annotation, I captured the information and use socket.io to sync all the canvases.' + $nl + '5. When people are using the knowledge graph, the socket.io can update the history dashboard to let everyone see the discussion and annotation.'

$ws.Range("C67").Value = '1. Dynamically update news in the main page.' + $nl + '2. Pass the Selected News information to the Start Chatting Room Page so that people can select both the News and the Room Number they want.' + $nl + '3. Using socket.io to sync the knowledge graph results to make the discussion more convenient.' + $nl + '4. Using Socket.io to note everyone the changes in the discussion dashboard.'

$ws.Range("B72").Value = '1. Implemented the socket.io to create and join a room.' + $nl + '2. Utilising socket.io to initialise the canvases in chatting room.' + $nl + '3. captured the drawing annotation information and using socket.io to update the canvases.' + $nl + '4. Implemented the knowledge graph and using socket.io to sync the searching result in chatting rooms.' + $nl + '5. Creating the canvas to take annotations on the searched pictures, and sync the annotation in rooms using socket.io.' + $nl + '6. Front-end Pages Design and implemented file selection to transform the picture into base64 to store the information in MongoDB.' + $nl + '7. Post the server to get news from the database, and dynamically update the returned data in the main page.' + $nl + '8. Onclick title or pictures to create a new chatting room, and pass the information to the Start Chatting Room Page.' + $nl + '9. Routes service configuration, parameters passing functions calling and render the correct pages.' + $nl + '10. Changes announcements and other notes when someone is joining the room, pushing new information or delete pictures in the chatting room. '

# Row height adjustments following the text edits above (content now wraps to more lines)
$ws.Rows.Item(56).RowHeight = 343.2
$ws.Rows.Item(67).RowHeight = 124.8

# Restore cursor/selection position as left by the author after editing
$ws.Range("E72").Select()
